$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 42605.88585648148
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B3").Value = -4
$ws.Range("C3").Value = 63
$ws.Range("D3").Value = 36
$ws.Range("E3").Value = 28
$ws.Range("F3").Value = 71
$ws.Range("G3").Value = 18298
$ws.Range("H3").Value = 4249
$ws.Range("I3").Value = 876
$ws.Range("J3").Value = 75
$ws.Range("K3").Value = 44
$ws.Range("L3").Value = 4
$ws.Range("M3").Value = 10
$ws.Range("N3").Value = "Noun"
